$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.779.67'
$ws.Range('E2').Value = '  -1.33%  '

$ws.Range('D3').Value = '3.856.85'
$ws.Range('E3').Value = '  +1.93%  '

$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  -0.33%  '

$ws.Range('D5').Value = "'424.08"
$ws.Range('E5').Value = '  +1.10%  '

$ws.Range('D6').Value = "'128.99"
$ws.Range('E6').Value = '  -2.61%  '

$ws.Range('D7').Value = '3.852.40'
$ws.Range('E7').Value = '  +2.28%  '

$ws.Range('D8').Value = "'0.605"
$ws.Range('E8').Value = '  -6.55%  '

$ws.Range('D9').Value = "'0.999"
$ws.Range('E9').Value = '  -0.05%  '

$ws.Range('D10').Value = "'0.720"
$ws.Range('E10').Value = '  -6.38%  '

$ws.Range('D11').Value = "'0.165"
$ws.Range('E11').Value = '  -10.86%  '

$ws.Range('D12').Value = "'0.0000356"
$ws.Range('E12').Value = '  -13.94%  '

$ws.Range('D13').Value = "'40.14"
$ws.Range('E13').Value = '  -5.95%  '

$ws.Range('D14').Value = '4.438.35'
$ws.Range('E14').Value = '  +1.69%  '

$ws.Range('D15').Value = "'10.01"
$ws.Range('E15').Value = '  -3.75%  '

$ws.Range('D16').Value = "'15.78"
$ws.Range('E16').Value = '  +19.27%  '

$ws.Range('D17').Value = '3.846.33'
$ws.Range('E17').Value = '  +2.16%  '

$ws.Range('E18').Value = '  -1.73%  '

$ws.Range('D19').Value = "'19.51"
$ws.Range('E19').Value = '  -4.86%  '

$ws.Range('D20').Value = '66.907.07'
$ws.Range('E20').Value = '  -1.20%  '

$ws.Range('D21').Value = "'1.07"
$ws.Range('E21').Value = '  -5.31%  '

$ws.Range('D22').Value = "'404.27"
$ws.Range('E22').Value = '  -9.33%  '

$ws.Range('D23').Value = "'14.25"
$ws.Range('E23').Value = '  -10.01%  '

$ws.Range('D24').Value = "'84.19"
$ws.Range('E24').Value = '  -7.40%  '

$ws.Range('D25').Value = "'2.97"
$ws.Range('E25').Value = '  -3.43%  '

$ws.Range('D26').Value = "'37.15"
$ws.Range('E26').Value = '  -2.83%  '

$ws.Range('D27').Value = "'5.83"
$ws.Range('E27').Value = '  +14.62%  '

$ws.Range('D28').Value = "'3.20"
$ws.Range('E28').Value = '  -3.88%  '

$ws.Range('D29').Value = "'9.42"
$ws.Range('E29').Value = '  -6.66%  '

$ws.Range('D30').Value = "'708.91"
$ws.Range('E30').Value = '  +3.30%  '

$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').Value = "'2.77"
$ws.Range('E31').Value = '  -0.07%  '

$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = "'0.121"
$ws.Range('E32').Value = '  -2.58%  '

$ws.Range('D33').Value = "'12.28"
$ws.Range('E33').Value = '  -2.88%  '

$ws.Range('D34').Value = "'7.48"
$ws.Range('E34').Value = '  +4.34%  '

$ws.Range('D35').Value = "'0.150"
$ws.Range('E35').Value = '  -8.97%  '

$ws.Range('B36').Value = 'InjectiveProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D36').Value = "'37.78"
$ws.Range('E36').Value = '  -8.95%  '

$ws.Range('B37').Value = 'Dai'
$ws.Range('C37').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D37').Value = "'1.00"
$ws.Range('E37').Value = '  +0.06%  '

$ws.Range('B38').Value = 'PEPE'
$ws.Range('C38').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D38').Value = '0.0₃0783'
$ws.Range('E38').Value = '  +7.09%  '

$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').Value = "'54.78"
$ws.Range('E39').Value = '  -5.17%  '

$ws.Range('D40').Value = "'0.0452"
$ws.Range('E40').Value = '  -7.69%  '

$ws.Range('D41').Value = "'2.92"
$ws.Range('E41').Value = '  -2.70%  '

$ws.Range('E42').Value = '  +0.18%  '

$ws.Range('D43').Value = "'0.135"
$ws.Range('E43').Value = '  -8.91%  '

$ws.Range('B44').Value = 'NEARProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D44').Value = "'4.42"
$ws.Range('E44').Value = '  +2.71%  '

$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D45').Value = "'3.15"
$ws.Range('E45').Value = '  -1.35%  '

$ws.Range('B46').Value = 'LidoDAOToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D46').Value = "'3.30"
$ws.Range('E46').Value = '  -2.97%  '

$ws.Range('D47').Value = "'143.74"
$ws.Range('E47').Value = '  -3.07%  '

$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').Value = "'2.05"
$ws.Range('E48').Value = '  -3.16%  '

$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = "'26.12"
$ws.Range('E49').Value = '  -7.71%  '

$ws.Range('E50').Value = '  -4.46%  '

$ws.Range('D51').Value = "'2.73"
$ws.Range('E51').Value = '  -5.43%  '
